# Add a new worksheet "N=200000" with the insertion-sort timing results,
# placed after the last existing sheet ("N=150000").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "N=200000"

# Move the freshly-added sheet to the end, after the last pre-existing
# sheet (capture the "last sheet" reference *after* the Add, since Add()
# inserts before the active sheet and shifts positions).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $lastSheet)

# Worksheet references resolve positionally, so re-fetch by name after the
# Move() call above moved the sheet's position within the collection.
$ws = $wb.Worksheets.Item("N=200000")

$ws.Range("A1").Value = "Execução"
$ws.Range("B1").Value = "Tempo (ms)"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "1375340.3971 ms"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "1428060.7851 ms"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "1385795.1810 ms"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "1428351.6641 ms"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "1430670.7091 ms"

$ws.Range("A7").Value = "Média"
$ws.Range("B7").Value = "1409643.7473 ms"

$ws.Range("A8").Value = "Desvio Padrão"
$ws.Range("B8").Value = "26817.8136 ms"
